$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the historical years to column A, rows 4-9
$years = @(2019, 2018, 2017, 2016, 2015, 2014)
$row = 4
foreach ($y in $years) {
    $ws.Cells.Item($row, 1).Value = $y
    $row++
}

# Move selection to A10 (the cell after the last filled row), with no frozen/scrolled topLeftCell
$ws.Range("A10").Select()
